$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("AI Engineer II - Blue Ring", "Blue Origin", "Seattle, WA, US USA", 16.7, "AI Engineer, RAG, FAISS, Pinecone, TensorFlow, PyTorch, Kinesis, Docker, Kubernetes, CI/CD", "2026-02-20", "https://www.indeed.com/viewjob?jk=ed2387d63f350d09"),
    @("AI Engineer III - Blue Ring", "Blue Origin", "Seattle, WA, US USA", 14.4, "AI Engineer, RAG, TensorFlow, PyTorch, Kinesis, Docker, Kubernetes, CI/CD, Databricks, Kafka", "2026-02-20", "https://www.indeed.com/viewjob?jk=282fd09ebf907b50"),
    @("AI Commercial & ML Ops Engineer", "MGM Resorts International", "Las Vegas, NV, US USA", 14.4, "Data Scientist, TensorFlow, PyTorch, MLflow, Docker, Kubernetes, CI/CD, Git, Databricks, PySpark", "2026-02-20", "https://www.indeed.com/viewjob?jk=1d1ac103ab878235"),
    @("Senior Data Engineer", "American Academy of Orthopaedic Surgeons", "Rosemont, IL, US USA", 13.3, "RAG, Synapse, Dataflow, CI/CD, Terraform, Git, Snowflake, PySpark, SQL, R", "2026-02-21", "https://www.indeed.com/viewjob?jk=c11ce56e07f4ef5b"),
    @("Senior Data Scientist", "Huntington Bank", "Detroit, MI, US USA", 11.1, "Data Scientist, RAG, TensorFlow, AWS SageMaker, NoSQL, Python, SQL, R, Scala, Bayesian", "2026-02-20", "https://www.indeed.com/viewjob?jk=05207a6563aa981a"),
    @("Senior Software Engineer", "eCapital", "US USA", 10, "RAG, Copilot, Docker, Kubernetes, Git, PostgreSQL, SQL, R, Scala", "2026-02-21", "https://www.indeed.com/viewjob?jk=98e272ab19ecdfca")
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).NumberFormat = "@"
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 6).Style = "Normal"
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}
